$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 12071.929
$ws.Range("I28").Value = 7789.778
$ws.Range("K28").Value = 7789.778
$ws.Range("M28").Value = -7304.778
$ws.Range("H40").Value = 1525.1875
$ws.Range("I40").Value = 1571.6428
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 1571.6428
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -1396.6428
$ws.Range("N40").Value = -1550
$ws.Range("H58").Value = 6957.375
$ws.Range("I58").Value = 8943.166999999999
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 26829.501
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -26679.501
$ws.Range("N58").Value = -3300
$ws.Range("H62").Value = 38098810
$ws.Range("J62").Value = 4992
$ws.Range("L62").Value = 4992
$ws.Range("N62").Value = -6240
$ws.Range("H65").Value = 38098810
$ws.Range("J65").Value = 4992
$ws.Range("L65").Value = 24960
$ws.Range("N65").Value = -31200
$ws.Range("H74").Value = 7678.4287
$ws.Range("I74").Value = 5937.25
$ws.Range("K74").Value = 5937.25
$ws.Range("M74").Value = -5001.25
$ws.Range("H77").Value = 7678.4287
$ws.Range("I77").Value = 5937.25
$ws.Range("K77").Value = 29686.25
$ws.Range("M77").Value = -25006.25
$ws.Range("H127").Value = 1482.95
$ws.Range("I127").Value = 880.1667
$ws.Range("K127").Value = 2640.5001
$ws.Range("M127").Value = 2319.4999
$ws.Range("H132").Value = 1811
$ws.Range("I132").Value = 1503.1428
$ws.Range("K132").Value = 4509.428400000001
$ws.Range("M132").Value = -1979.428400000001
$ws.Range("H135").Value = 8338471.5
$ws.Range("I135").Value = 25002030
$ws.Range("J135").Value = 6691.875
$ws.Range("K135").Value = 225018270
$ws.Range("L135").Value = 60226.875
$ws.Range("M135").Value = -225015735
$ws.Range("N135").Value = -65296.875
$ws.Range("H137").Value = 9499.272000000001
$ws.Range("I137").Value = 1083.875
$ws.Range("J137").Value = 31940.334
$ws.Range("K137").Value = 3251.625
$ws.Range("L137").Value = 95821.00199999999
$ws.Range("M137").Value = -701.625
$ws.Range("N137").Value = -100921.002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8265727
$ws.Range("I2").Value = 11364625
$ws.Range("K2").Value = 11364625
$ws.Range("M2").Value = -11364512
$ws.Range("H42").Value = 22999.5
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30972
$ws.Range("H44").Value = 92500
$ws.Range("J44").Value = 92500
$ws.Range("L44").Value = 92500
$ws.Range("N44").Value = -93476
$ws.Range("H45").Value = 2628.2
$ws.Range("I45").Value = 2035.3334
$ws.Range("K45").Value = 2035.3334
$ws.Range("M45").Value = -1658.3334
$ws.Range("H114").Value = 35000
$ws.Range("J114").Value = 35000
$ws.Range("L114").Value = 35000
$ws.Range("N114").Value = -43678
$ws.Range("H116").Value = 8265727
$ws.Range("I116").Value = 11364625
$ws.Range("K116").Value = 11364625
$ws.Range("M116").Value = -11362331

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8265727
$ws.Range("I3").Value = 11364625
$ws.Range("K3").Value = 11364625
$ws.Range("M3").Value = -11364511
$ws.Range("H86").Value = 138028.4
$ws.Range("I86").Value = 5124.6
$ws.Range("K86").Value = 5124.6
$ws.Range("M86").Value = -4001.6
$ws.Range("H89").Value = 138028.4
$ws.Range("I89").Value = 5124.6
$ws.Range("K89").Value = 25623
$ws.Range("M89").Value = -20007

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3888.8518
$ws.Range("I31").Value = 1555.4615
$ws.Range("J31").Value = 6055.5713
$ws.Range("K31").Value = 1555.4615
$ws.Range("L31").Value = 6055.5713
$ws.Range("M31").Value = -1260.4615
$ws.Range("N31").Value = -6645.5713
$ws.Range("H34").Value = 3888.8518
$ws.Range("I34").Value = 1555.4615
$ws.Range("J34").Value = 6055.5713
$ws.Range("K34").Value = 1555.4615
$ws.Range("L34").Value = 6055.5713
$ws.Range("M34").Value = -1353.4615
$ws.Range("N34").Value = -6459.5713
$ws.Range("H58").Value = 4164.9653
$ws.Range("J58").Value = 6990.8184
$ws.Range("L58").Value = 6990.8184
$ws.Range("N58").Value = -7396.8184
$ws.Range("H99").Value = 1637.3334
$ws.Range("I99").Value = 1764.8
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1764.8
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -266.8
$ws.Range("N99").Value = -3996
$ws.Range("H126").Value = 1637.3334
$ws.Range("I126").Value = 1764.8
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 5294.4
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -2824.4
$ws.Range("N126").Value = -7940
$ws.Range("H136").Value = 4164.9653
$ws.Range("J136").Value = 6990.8184
$ws.Range("L136").Value = 20972.4552
$ws.Range("N136").Value = -26072.4552

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2600.625
$ws.Range("J34").Value = 3325.9167
$ws.Range("L34").Value = 9977.750100000001
$ws.Range("N34").Value = -10145.7501
$ws.Range("H55").Value = 933
$ws.Range("I55").Value = 253
$ws.Range("J55").Value = 1499.6666
$ws.Range("K55").Value = 759
$ws.Range("L55").Value = 4498.9998
$ws.Range("M55").Value = -582
$ws.Range("N55").Value = -4852.9998
$ws.Range("H57").Value = 2000
$ws.Range("I57").Value = 1500
$ws.Range("K57").Value = 4500
$ws.Range("M57").Value = -3941
$ws.Range("H74").Value = 10428.571
$ws.Range("J74").Value = 10664.5
$ws.Range("L74").Value = 31993.5
$ws.Range("N74").Value = -34115.5
$ws.Range("H77").Value = 10428.571
$ws.Range("J77").Value = 10664.5
$ws.Range("L77").Value = 95980.5
$ws.Range("N77").Value = -106588.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16554.889
$ws.Range("I70").Value = 12499.8
$ws.Range("J70").Value = 21623.75
$ws.Range("K70").Value = 12499.8
$ws.Range("L70").Value = 21623.75
$ws.Range("M70").Value = -12229.8
$ws.Range("N70").Value = -22163.75
$ws.Range("H73").Value = 16554.889
$ws.Range("I73").Value = 12499.8
$ws.Range("J73").Value = 21623.75
$ws.Range("K73").Value = 12499.8
$ws.Range("L73").Value = 21623.75
$ws.Range("M73").Value = -11563.8
$ws.Range("N73").Value = -23495.75
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -56134
$ws.Range("H132").Value = 3780.65
$ws.Range("I132").Value = 2115.4
$ws.Range("J132").Value = 8776.4
$ws.Range("K132").Value = 6346.200000000001
$ws.Range("L132").Value = 26329.2
$ws.Range("M132").Value = -3816.200000000001
$ws.Range("N132").Value = -31389.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12721.066
$ws.Range("I40").Value = 13293.538
$ws.Range("K40").Value = 13293.538
$ws.Range("M40").Value = -13157.538
$ws.Range("H122").Value = 4321.636
$ws.Range("I122").Value = 3375.6155
$ws.Range("K122").Value = 10126.8465
$ws.Range("M122").Value = -7676.8465
$ws.Range("H132").Value = 9066.375
$ws.Range("I132").Value = 8449.700000000001
$ws.Range("J132").Value = 10094.167
$ws.Range("K132").Value = 25349.1
$ws.Range("L132").Value = 30282.501
$ws.Range("M132").Value = -22819.1
$ws.Range("N132").Value = -35342.501

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("H107").Value = 2165.889
$ws.Range("I107").Value = 2164.7693
$ws.Range("K107").Value = 6494.3079
$ws.Range("M107").Value = -4574.3079
$ws.Range("H122").Value = 4496.579
$ws.Range("I122").Value = 4628.2
$ws.Range("K122").Value = 13884.6
$ws.Range("M122").Value = -11434.6
$ws.Range("H126").Value = 3161.5715
$ws.Range("I126").Value = 3161.5715
$ws.Range("K126").Value = 9484.7145
$ws.Range("M126").Value = -7014.7145
$ws.Range("H132").Value = 5676.879
$ws.Range("I132").Value = 2856.3684
$ws.Range("K132").Value = 8569.1052
$ws.Range("M132").Value = -6039.1052
$ws.Range("N101").ClearContents()
